$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new cell value (adds "mango women" to the shared strings table)
$ws.Range("C1").Value = "mango women"

# Update the active cell / selection to A8
$ws.Range("A8").Select()
